# ------------------------------------------------------------------
# GDD.docx edit: add "General Information" / "External influences"
# sections, split "Gameplay" into Controls/Mechanics/Victory sections,
# drop the old spell-check proofErr markers, merge some runs, and tag
# every paragraph/run with an explicit en-US language.
# ------------------------------------------------------------------

$d = $word.ActiveDocument

# --------------------------------------------------------------
# 1. Create the "List Paragraph" style (Italian template styleId
#    "Paragrafoelenco") used by the new bullet list.
# --------------------------------------------------------------
$paraStyle = $d.Styles.Add("Paragrafoelenco", 1)
$paraStyle.BaseStyle = $d.Styles("Normale")
$paraStyle.NameLocal = "List Paragraph"
$paraStyle.Priority = 34
$paraStyle.QuickStyle = $true
$paraStyle.ParagraphFormat.LeftIndent = 36
$paraStyle.NoSpaceBetweenParagraphsOfSameStyle = $true

# --------------------------------------------------------------
# 2. Bootstrap a numbering definition (numId 1 / abstractNumId 0,
#    the standard Word bullet template) by applying it to a
#    throw-away paragraph, then remove that paragraph again. The
#    numbering.xml part (and its abstractNum/num) stays behind.
# --------------------------------------------------------------
$bootstrap = $d.Paragraphs.Add()
$bootstrap.Range.Text = "x"
$bootstrap.Range.ListFormat.ApplyListTemplateWithLevel($word.ListGalleries.Item(1).ListTemplates.Item(1))
$bootstrap.Range.Delete()

# --------------------------------------------------------------
# 3. Clear out all existing paragraphs, leaving exactly one empty
#    paragraph (so the section properties / sectPr stay intact).
# --------------------------------------------------------------
while ($d.Paragraphs.Count -gt 1) {
    $d.Paragraphs.Item(1).Range.Delete()
}
$lastPara = $d.Paragraphs.Item(1)
$d.Range($lastPara.Range.Start, $lastPara.Range.End - 1).Delete()

# --------------------------------------------------------------
# 4. Insert the full set of new paragraphs (everything except the
#    very last, empty, paragraph) as raw OOXML at the top of the
#    document, ahead of the single leftover empty paragraph.
# --------------------------------------------------------------
$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$lang = '<w:rPr><w:lang w:val="en-US"/></w:rPr>'

$body = @"
<w:p $w>
  <w:pPr><w:pStyle w:val="Titolo1"/>$lang</w:pPr>
  <w:r>$lang<w:t>GDD</w:t></w:r>
</w:p>
<w:p $w>
  <w:pPr><w:pStyle w:val="Titolo2"/>$lang</w:pPr>
  <w:r>$lang<w:t>General Information</w:t></w:r>
</w:p>
<w:p $w>
  <w:pPr>
    <w:pStyle w:val="Paragrafoelenco"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
    $lang
  </w:pPr>
  <w:r>$lang<w:t>Single player</w:t></w:r>
</w:p>
<w:p $w>
  <w:pPr>
    <w:pStyle w:val="Paragrafoelenco"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
    $lang
  </w:pPr>
  <w:r>$lang<w:t>Genre: fps, platform, puzzle</w:t></w:r>
</w:p>
<w:p $w>
  <w:pPr>
    <w:pStyle w:val="Paragrafoelenco"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
    $lang
  </w:pPr>
  <w:r>$lang<w:t>Platform: PC</w:t></w:r>
</w:p>
<w:p $w>
  <w:pPr><w:pStyle w:val="Titolo2"/>$lang</w:pPr>
  <w:r>$lang<w:t>External influences</w:t></w:r>
</w:p>
<w:p $w>
  <w:pPr>$lang</w:pPr>
</w:p>
<w:p $w>
  <w:pPr><w:pStyle w:val="Titolo2"/>$lang</w:pPr>
  <w:r>$lang<w:t>Story</w:t></w:r>
  <w:r>$lang<w:t xml:space="preserve"> and</w:t></w:r>
  <w:r>$lang<w:t xml:space="preserve"> setting</w:t></w:r>
</w:p>
<w:p $w>
  <w:pPr>$lang</w:pPr>
</w:p>
<w:p $w>
  <w:pPr><w:pStyle w:val="Titolo2"/>$lang</w:pPr>
  <w:r>$lang<w:t>Gameplay</w:t></w:r>
</w:p>
<w:p $w>
  <w:pPr><w:pStyle w:val="Titolo3"/>$lang</w:pPr>
  <w:r>$lang<w:t>Controls</w:t></w:r>
</w:p>
<w:p $w>
  <w:pPr>$lang</w:pPr>
</w:p>
<w:p $w>
  <w:pPr><w:pStyle w:val="Titolo3"/>$lang</w:pPr>
  <w:r>$lang<w:t>Mechanics</w:t></w:r>
</w:p>
<w:p $w>
  <w:pPr>$lang</w:pPr>
  <w:r>$lang<w:t>First person view, WASD movement</w:t></w:r>
  <w:r>$lang<w:t>, jumping</w:t></w:r>
  <w:r>$lang<w:t>, shooting</w:t></w:r>
  <w:r>$lang<w:t>, enemies</w:t></w:r>
  <w:r>$lang<w:t>’ AI.</w:t></w:r>
</w:p>
<w:p $w>
  <w:pPr><w:pStyle w:val="Titolo3"/>$lang</w:pPr>
  <w:r>$lang<w:t>Victory and Game Goals</w:t></w:r>
</w:p>
"@

$insertionPoint = $d.Range(0, 0)
$insertionPoint.InsertXML($body)

# --------------------------------------------------------------
# 5. The single leftover empty paragraph is now the very last
#    paragraph in the document; just tag it with en-US too.
# --------------------------------------------------------------
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalPara.Range.LanguageID = "en-US"

Write-Output "done"
